$wb = $excel.ActiveWorkbook

# --- Docentes sheet: rename header "estadia" -> "horas extra clase",
#     widen column H to fit the new text, and make it the active/selected sheet.
$wsDocentes = $wb.Worksheets.Item("Docentes")
$wsDocentes.Range("H1").Value = "horas extra clase"
$wsDocentes.Range("H1").EntireColumn.ColumnWidth = 14.33

# --- Grupos sheet was previously the active tab; Docentes becomes active instead.
# Activating Docentes flips Docentes.tabSelected = true and Grupos.tabSelected =
# false automatically, and moves the workbook's activeTab to sheet index 0.
$wsDocentes.Activate()
$wsDocentes.Range("I5").Select()
